$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 1999.5
$ws.Range("I26").Value = 1999.5
$ws.Range("K26").Value = 1999.5
$ws.Range("M26").Value = -1655.5
# Row 113
$ws.Range("H113").Value = 1797.8
$ws.Range("I113").Value = 2001
$ws.Range("J113").Value = 985
$ws.Range("K113").Value = 2001
$ws.Range("L113").Value = 985
$ws.Range("M113").Value = 1253
$ws.Range("N113").Value = -7493
# Row 125
$ws.Range("H125").Value = 1135.4286
$ws.Range("J125").Value = 1031
$ws.Range("L125").Value = 9279
$ws.Range("N125").Value = -14199
# Row 132
$ws.Range("H132").Value = 3125.2856
$ws.Range("I132").Value = 1175.8
$ws.Range("K132").Value = 3527.4
$ws.Range("M132").Value = -997.3999999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 50
$ws.Range("H50").Value = 20509.5
$ws.Range("I50").Value = 12697
$ws.Range("J50").Value = 25197
$ws.Range("K50").Value = 12697
$ws.Range("L50").Value = 25197
$ws.Range("M50").Value = -11983
$ws.Range("N50").Value = -26625
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 61
$ws.Range("H61").Value = 2128.3333
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 2385
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 2385
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -2809
# Row 94
$ws.Range("H94").Value = 56666
$ws.Range("J94").Value = 56666
$ws.Range("L94").Value = 56666
$ws.Range("N94").Value = -58468
# Row 136
$ws.Range("H136").Value = 2128.3333
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2385
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 7155
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -12255

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 33
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -4664
$ws.Range("N33").Value = -5672
# Row 36
$ws.Range("H36").Value = 8650
$ws.Range("I36").Value = 9187.5
$ws.Range("J36").Value = 6500
$ws.Range("K36").Value = 9187.5
$ws.Range("L36").Value = 6500
$ws.Range("M36").Value = -8653.5
$ws.Range("N36").Value = -7568
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 26
$ws.Range("H26").Value = 2000
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1713
# Row 36
$ws.Range("H36").Value = 5666.3335
$ws.Range("I36").Value = 5666.3335
$ws.Range("K36").Value = 5666.3335
$ws.Range("M36").Value = -5278.3335
# Row 38
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -9623
$ws.Range("N38").ClearContents()
# Row 40
$ws.Range("H40").Value = 5666.3335
$ws.Range("I40").Value = 5666.3335
$ws.Range("K40").Value = 5666.3335
$ws.Range("M40").Value = -5506.3335
# Row 42
$ws.Range("H42").Value = 19000
$ws.Range("I42").Value = 6000
$ws.Range("K42").Value = 6000
$ws.Range("M42").Value = -5407
# Row 44
$ws.Range("H44").Value = 25999.8
$ws.Range("I44").Value = 23333
$ws.Range("K44").Value = 23333
$ws.Range("M44").Value = -22891
# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -9789
$ws.Range("N46").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 1099.8
$ws.Range("J81").Value = 500
$ws.Range("L81").Value = 1500
$ws.Range("N81").Value = -3746
# Row 84
$ws.Range("H84").Value = 1099.8
$ws.Range("J84").Value = 500
$ws.Range("L84").Value = 4500
$ws.Range("N84").Value = -15732

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 195
$ws.Range("I31").Value = 195
$ws.Range("K31").Value = 195
$ws.Range("M31").Value = 97
# Row 37
$ws.Range("H37").Value = 195
$ws.Range("I37").Value = 195
$ws.Range("K37").Value = 195
$ws.Range("M37").Value = 82
# Row 49
$ws.Range("H49").Value = 5807.6924
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 10250
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 10250
$ws.Range("M49").Value = -4816
$ws.Range("N49").Value = -10618
# Row 54
$ws.Range("H54").Value = 10086
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 93
$ws.Range("H93").Value = 89999.336
$ws.Range("I93").Value = 90000
$ws.Range("J93").Value = 89999
$ws.Range("K93").Value = 90000
$ws.Range("L93").Value = 89999
$ws.Range("M93").Value = -88128
$ws.Range("N93").Value = -93743

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Range("H24").Value = 4000
$ws.Range("J24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("N24").Value = -4686
# Row 26
$ws.Range("H26").Value = 3928.4285
$ws.Range("I26").Value = 1249.5
$ws.Range("K26").Value = 1249.5
$ws.Range("M26").Value = -954.5
# Row 35
$ws.Range("H35").Value = 2513.75
$ws.Range("I35").Value = 1385
$ws.Range("J35").Value = 5900
$ws.Range("K35").Value = 1385
$ws.Range("L35").Value = 5900
$ws.Range("M35").Value = -1049
$ws.Range("N35").Value = -6572
# Row 39
$ws.Range("H39").Value = 2500
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 1000
$ws.Range("M39").Value = -540
# Row 40
$ws.Range("H40").Value = 41357.43
$ws.Range("I40").Value = 100000
$ws.Range("K40").Value = 100000
$ws.Range("M40").Value = -99864

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 58407.332
$ws.Range("J28").Value = 58407.332
$ws.Range("L28").Value = 58407.332
$ws.Range("N28").Value = -59103.332
# Row 29
$ws.Range("H29").Value = 39260
$ws.Range("I29").Value = 29075
$ws.Range("K29").Value = 29075
$ws.Range("M29").Value = -28785
# Row 32
$ws.Range("H32").Value = 3463
$ws.Range("I32").Value = 3463
$ws.Range("K32").Value = 3463
$ws.Range("M32").Value = -3146
# Row 70
$ws.Range("H70").Value = 90000
$ws.Range("I70").Value = 90000
$ws.Range("K70").Value = 90000
$ws.Range("M70").Value = -89685
# Row 73
$ws.Range("H73").Value = 90000
$ws.Range("I73").Value = 90000
$ws.Range("K73").Value = 90000
$ws.Range("M73").Value = -88908
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
